# A-O_Demand_COMPLETED.xlsx
#
# The "nan" placeholder text that used to fill column AP (the last column,
# "Target.Unit") for most data rows is cleared out entirely, while the rows
# that legitimately carry a unit ("Gvkm") keep that value.
#
# Column AP = 42nd column. Row 1 is the header row and is left untouched.
# Rows 74, 75, 76, 78, 79, 80, 84, 86 and 89 keep "Gvkm" in AP; every other
# data row (2-90) has its AP cell cleared completely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$apCol = 42
$keepGvkmRows = @(74, 75, 76, 78, 79, 80, 84, 86, 89)

for ($r = 2; $r -le 90; $r++) {
    if ($keepGvkmRows -contains $r) {
        $ws.Cells.Item($r, $apCol).Value = "Gvkm"
    } else {
        $ws.Cells.Item($r, $apCol).ClearContents()
    }
}
